$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Etat Taxes")

# --- Row 2 ---
$ws.Range("A2").Value = "004/ZZZ/AV2"
$ws.Range("C2").Value = "IR801997"
$ws.Range("D2").Value = "NOUBAIL MOHAMMED"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "--"
$ws.Range("J2").Value = "--"
$ws.Range("L2").Value = 3000
$ws.Range("O2").Value = 3000

# --- Row 3 ---
$ws.Range("A3").Value = "004/ZZZ/AV2"
$ws.Range("C3").Value = "Q251990"
$ws.Range("D3").Value = "NOUBAIL MOUNTASSIR"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "--"
$ws.Range("J3").Value = "--"
$ws.Range("L3").Value = 3000
$ws.Range("O3").Value = 3000

# --- Row 4 ---
$ws.Range("A4").Value = "004/ZZZ/AV2"
$ws.Range("C4").Value = "IR801997"
$ws.Range("D4").Value = "NOUBAIL MOHAMMED"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 0
$ws.Range("O4").Value = 1000

# --- Row 5 ---
$ws.Range("A5").Value = "004/ZZZ/AV2"
$ws.Range("B5").Value = "Direction régionale"
$ws.Range("C5").Value = "Q251990"
$ws.Range("D5").Value = "NOUBAIL MOUNTASSIR"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1000
$ws.Range("J5").Value = 0
$ws.Range("O5").Value = 1000

# --- Row 6 (totals) ---
$ws.Range("H6").Value = 2000
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 6000
$ws.Range("O6").Value = 8000
